$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("AddCustomerTest")

# --- Update existing sheet (AddCustomerTest) - part 1 ---
$ws1.Range("D1").Value = "alerttext"
$ws1.Range("C2").Value = 411033
$ws1.Range("D2").Value = "Customer added successfully"

# --- Add the new sheet right after AddCustomerTest ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "OpenAccountTest"

$ws2.Range("A1").Value = "customer"
$ws2.Range("B1").Value = "currency"
$ws2.Range("A2").Value = "sandip thopate"
$ws2.Range("B2").Value = "rupee"

$ws2.Columns.Item(1).ColumnWidth = 15.5
$ws2.Columns.Item(2).ColumnWidth = 12.75

# --- Finish updating AddCustomerTest (row 3) ---
$ws1.Range("A3").Value = "vishal"
$ws1.Range("B3").Value = "sahu"
$ws1.Range("C3").Value = 411034
$ws1.Range("D3").Value = "Customer added successfully"

# --- Make OpenAccountTest the active sheet/tab with H17 selected ---
$ws2.Range("H17").Select()
